$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129, shifting existing rows 129..177 down to 130..178
$ws.Rows("129").Insert()

# Populate the newly inserted row 129 with its data
$ws.Range("A129").Value = 3
$ws.Range("B129").Value = "Femacal de La Calera"
$ws.Range("C129").Value = "Coquimbo"
$ws.Range("D129").Value = 44784
$ws.Range("E129").Value = 5
$ws.Range("F129").Value = 100112026
$ws.Range("G129").Value = "Haba"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 83
$ws.Range("K129").Value = 15000
$ws.Range("L129").Value = 16000
$ws.Range("M129").Value = 15458
$ws.Range("N129").Value = "$/saco 25 kilos"
$ws.Range("O129").Value = "Provincia de Limar" + [char]0x00ED
$ws.Range("P129").Value = 618
$ws.Range("Q129").Value = 25
$ws.Range("R129").Value = "Hortaliza"
